$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: add the new 2022-Q4 row and shift the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy the style (bold/centered/bordered) of the existing A2 index cell onto
# the brand-new A6 cell before writing into it, so it matches the look of
# the other index cells (A2:A5) instead of staying plain.
$summary.Range("A2").Copy()
$summary.Range("A6").PasteSpecial(-4122)

$summaryData = @(
    @(0, "2022-Q4", 15, 3.14),
    @(1, "2022-Q3", 19, 3.66),
    @(2, "2022-Q2", 9, 2.32),
    @(3, "2022-Q1", 8, 0.39),
    @(4, "2021-Q4", 4, 2.64)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $r = $i + 2
    $row = $summaryData[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that used to be "2022-Q3"), holding the quarter's fund holdings.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Reuse the header/index styling from the summary sheet instead of inventing
# new style combinations: B1 carries the bold/centered/bordered "header"
# style used across the workbook, and A2 carries the same style used for the
# numeric index column.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A16").PasteSpecial(-4122)

# Columns B:G hold values that look numeric (fund codes with leading zeros,
# decimal percentages, etc.) but must stay as literal text, exactly like the
# equivalent columns on the other quarter sheets. Formatting the range as
# Text before typing the values stops Excel from re-interpreting them as
# numbers (and stripping the leading zeros / trailing precision).
$q4.Range("B2:G16").NumberFormat = "@"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$funds = @(
    @("050001", "博时价值增长混合",             "21.87", "74.33", "3.12", "0.6823", 5),
    @("014038", "交银启诚混合A",                 "24.58", "81.04", "2.66", "0.6538", 6),
    @("001128", "宝盈新兴产业灵活配置混合A",     "9.13",  "94.44", "4.47", "0.4081", 3),
    @("001877", "宝盈国家安全沪港深股票A",       "6.58",  "94.32", "4.74", "0.3119", 3),
    @("050201", "博时价值增长贰号混合",           "9.54",  "74.96", "3.17", "0.3024", 6),
    @("001487", "宝盈优势产业灵活配置混合A",     "10.62", "94.48", "2.73", "0.2899", 9),
    @("014039", "交银启诚混合C",                 "7.22",  "81.04", "2.66", "0.1921", 6),
    @("001075", "宝盈转型动力灵活配置混合A",     "4.42",  "91.93", "3.44", "0.1520", 6),
    @("012771", "宝盈优势产业灵活配置混合C",     "3.19",  "94.48", "2.73", "0.0871", 9),
    @("012815", "宝盈新兴产业灵活配置混合C",     "1.06",  "94.44", "4.47", "0.0474", 3),
    @("013613", "宝盈国家安全沪港深股票C",       "0.13",  "94.32", "4.74", "0.0062", 3),
    @("007316", "交银施罗德可转债债券A",         "0.82",  "20.54", "0.67", "0.0055", 6),
    @("007317", "交银施罗德可转债债券C",         "0.22",  "20.54", "0.67", "0.0015", 6),
    @("001978", "泰信互联网+主题灵活配置混合",   "0.06",  "77.79", "2.11", "0.0013", 7),
    @("015389", "宝盈转型动力灵活配置混合C",     "0.03",  "91.93", "3.44", "0.0010", 6)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $r = $i + 2
    $f = $funds[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $f[0]
    $q4.Cells.Item($r, 3).Value = $f[1]
    $q4.Cells.Item($r, 4).Value = $f[2]
    $q4.Cells.Item($r, 5).Value = $f[3]
    $q4.Cells.Item($r, 6).Value = $f[4]
    $q4.Cells.Item($r, 7).Value = $f[5]
    $q4.Cells.Item($r, 8).Value = $f[6]
}
